# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to the Titan_Profits leve-profit workbook
# (columns H..N = price/profit calculations) per sheet/row, as described
# by the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 53: No Accounting for Waste
$ws.Range("H53").Value = 404.42856
$ws.Range("I53").Value = 417.85715
$ws.Range("J53").Value = 391
$ws.Range("K53").Value = 417.85715
$ws.Range("L53").Value = 391
$ws.Range("M53").Value = 219.14285
$ws.Range("N53").Value = -1665

# ALC row 62: The Mustache Suits Him
$ws.Range("H62").Value = 6723.2354
$ws.Range("I62").Value = 5109.5
$ws.Range("J62").Value = 9028.571
$ws.Range("K62").Value = 5109.5
$ws.Range("L62").Value = 9028.571
$ws.Range("M62").Value = -4485.5
$ws.Range("N62").Value = -10276.571

# ALC row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 6723.2354
$ws.Range("I65").Value = 5109.5
$ws.Range("J65").Value = 9028.571
$ws.Range("K65").Value = 25547.5
$ws.Range("L65").Value = 45142.855
$ws.Range("M65").Value = -22427.5
$ws.Range("N65").Value = -51382.855

# ALC row 100: Asking for a Friend
$ws.Range("H100").Value = 33335522
$ws.Range("I100").Value = 5000
$ws.Range("J100").Value = 41668150
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 41668150
$ws.Range("M100").Value = -4459
$ws.Range("N100").Value = -41669232

# ALC row 107: Another Man's Ink
$ws.Range("H107").Value = 508424.03
$ws.Range("I107").Value = 658798
$ws.Range("J107").Value = 911.875
$ws.Range("K107").Value = 658798
$ws.Range("L107").Value = 911.875
$ws.Range("M107").Value = -656878
$ws.Range("N107").Value = -4751.875

# ALC row 113: Amaro Kart
$ws.Range("H113").Value = 334668.34

# ALC row 133: Big Brush, Big Dreams
$ws.Range("H133").Value = 41666.668
$ws.Range("J133").Value = 41666.668
$ws.Range("L133").Value = 41666.668
$ws.Range("N133").Value = -51786.668

# ALC row 134: Binding Spells
$ws.Range("H134").Value = 50488
$ws.Range("J134").Value = 50488
$ws.Range("L134").Value = 50488
$ws.Range("N134").Value = -60628

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust
$ws.Range("H32").Value = 2076.6057
$ws.Range("I32").Value = 1860.5636
$ws.Range("K32").Value = 1860.5636
$ws.Range("M32").Value = -1573.5636

# ARM row 122: Haste for High Durium
$ws.Range("H122").Value = 2138.3333
$ws.Range("I122").Value = 1816.3334
$ws.Range("K122").Value = 5449.0002
$ws.Range("M122").Value = -2999.0002

$ws = $wb.Worksheets.Item("BSM")
# BSM row 64: With Bearings Straight
$ws.Range("H64").Value = 811.6
$ws.Range("J64").Value = 790
$ws.Range("L64").Value = 790
$ws.Range("N64").Value = -1240

# BSM row 67: Bearing the Brunt (L)
$ws.Range("H67").Value = 811.6
$ws.Range("J67").Value = 790
$ws.Range("L67").Value = 790
$ws.Range("N67").Value = -2350

# BSM row 99: Meddle in Metal
$ws.Range("H99").Value = 3508.182
$ws.Range("I99").Value = 1532.2222
$ws.Range("J99").Value = 12400
$ws.Range("K99").Value = 1532.2222
$ws.Range("L99").Value = 12400
$ws.Range("M99").Value = -34.22219999999993
$ws.Range("N99").Value = -15396

$ws = $wb.Worksheets.Item("CRP")
# CRP row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1442.6
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1442.6
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4327.799999999999
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9227.799999999999

# CRP row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1874.2333
$ws.Range("I134").Value = 728.5417
$ws.Range("K134").Value = 2185.6251
$ws.Range("M134").Value = 349.3748999999998

$ws = $wb.Worksheets.Item("CUL")
# CUL row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 57359.223
$ws.Range("I3").Value = 52025
$ws.Range("K3").Value = 156075
$ws.Range("M3").Value = -155963

# CUL row 5: What a Sap
$ws.Range("H5").Value = 1429.75
$ws.Range("I5").Value = 731.6667
$ws.Range("J5").Value = 1728.9286
$ws.Range("K5").Value = 2195.0001
$ws.Range("L5").Value = 5186.7858
$ws.Range("M5").Value = -2083.0001
$ws.Range("N5").Value = -5410.7858

# CUL row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 1506.45
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# CUL row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1429.75
$ws.Range("I135").Value = 731.6667
$ws.Range("J135").Value = 1728.9286
$ws.Range("K135").Value = 6585.0003
$ws.Range("L135").Value = 15560.3574
$ws.Range("M135").Value = -4050.0003
$ws.Range("N135").Value = -20630.3574

# CUL row 137: Creative Chocolate
$ws.Range("H137").Value = 4044403.2
$ws.Range("J137").Value = 131188.75
$ws.Range("L137").Value = 393566.25
$ws.Range("N137").Value = -403766.25

# CUL row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 1010.1667
$ws.Range("I138").Value = 1010.1667
$ws.Range("K138").Value = 3030.5001
$ws.Range("M138").Value = 2109.4999

# CUL row 139: Najoothie
$ws.Range("H139").Value = 1935.3478
$ws.Range("I139").Value = 1659.381
$ws.Range("J139").Value = 4833
$ws.Range("K139").Value = 4978.143
$ws.Range("L139").Value = 14499
$ws.Range("M139").Value = 161.857
$ws.Range("N139").Value = -24779

# CUL row 141: Ocean Explosion
$ws.Range("H141").Value = 7075.5557
$ws.Range("J141").Value = 7000
$ws.Range("L141").Value = 21000
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1589616.8
$ws.Range("I122").Value = 3705039.2
$ws.Range("K122").Value = 11115117.6
$ws.Range("M122").Value = -11112667.6

$ws = $wb.Worksheets.Item("LTW")
# LTW row 87: Bar of the Bannermen
$ws.Range("H87").Value = 285094.5
$ws.Range("J87").Value = 285094.5
$ws.Range("L87").Value = 285094.5
$ws.Range("N87").Value = -287340.5

# LTW row 90: Do My Little Turn on the Stonewalk (L)
$ws.Range("H90").Value = 285094.5
$ws.Range("J90").Value = 285094.5
$ws.Range("L90").Value = 855283.5
$ws.Range("N90").Value = -866515.5

# LTW row 122: Hell on Leather
$ws.Range("H122").Value = 3189.2856
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3419.0476
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 10257.1428
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -15157.1428

# LTW row 132: Tenets of Tanning
$ws.Range("H132").Value = 3130.2307
$ws.Range("I132").Value = 2270.6897
$ws.Range("K132").Value = 6812.0691
$ws.Range("M132").Value = -4282.0691

$ws = $wb.Worksheets.Item("WVR")
# WVR row 99: Say Yes to Formal Dress
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# WVR row 122: Heavy Armoire
$ws.Range("H122").Value = 40361
$ws.Range("I122").Value = 85023.414
$ws.Range("J122").Value = 2078.9285
$ws.Range("K122").Value = 255070.242
$ws.Range("L122").Value = 6236.7855
$ws.Range("M122").Value = -252620.242
$ws.Range("N122").Value = -11136.7855

# WVR row 123: Helping Handwear
$ws.Range("H123").Value = 33257.4
$ws.Range("J123").Value = 33257.4
$ws.Range("L123").Value = 33257.4
$ws.Range("N123").Value = -43057.4

# WVR row 124: Hot Heads
$ws.Range("H124").Value = 45214.5
$ws.Range("J124").Value = 45214.5
$ws.Range("L124").Value = 45214.5
$ws.Range("N124").Value = -55034.5

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 12822747
$ws.Range("I132").Value = 17858880
$ws.Range("K132").Value = 53576640
$ws.Range("M132").Value = -53574110
